# ===== Edit script: MINAGRIS_Plastic.R_ReadMe.xlsx =====
# Commit: New plotting of the Bar graphs: Average per Field
$wb = $excel.ActiveWorkbook

# --- 'Merged WUR Ubern data': move selection to C5 ---
$ws3 = $wb.Worksheets.Item("Merged WUR Ubern data")
$ws3.Activate()
$ws3.Range("C5").Select()

# --- 'Summaries': insert row for new ID 8 (Farm), extend factor-legend with h/i ---
$ws4 = $wb.Worksheets.Item("Summaries")
$ws4.Rows.Item(9).Insert()
$ws4.Range("A9").Value = 8
$ws4.Range("B9").Value = 'Farm'

$ws4.Range("A18").Value = 'h'
$ws4.Range("C18").Value = 'Polymer.red3 * Size_cat.um'
$ws4.Range("A19").Value = 'i'
$ws4.Range("C19").Value = 'Polymer.red3'

# --- New sheet: 'Replicated soils in IR files ' ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $lastSheet)
$ws5.Name = 'Replicated soils in IR files '

# Row 4: sample-ID headers
$ws5.Range("C4").Value = '3.6.1_S1'
$ws5.Range("D4").Value = '6.1.1_S1'
$ws5.Range("E4").Value = '8.1.2_S1'
$ws5.Range("F4").Value = '8.3.1_S1'
$ws5.Range("G4").Value = '8.3.2_S1'
$ws5.Range("H4").Value = '8.5.2_S1'
$ws5.Range("I4").Value = '8.6.1_S1'
$ws5.Range("J4").Value = '8.6.2_S1'
$ws5.Range("K4").Value = '8.7.2_S1'
$ws5.Range("L4").Value = '8.8.1_S1'
$ws5.Range("M4").Value = '9.4.2_S1'
$ws5.Range("N4").Value = '11.10.2_S1'
$ws5.Range("O4").Value = '11.10.1_S1'
$ws5.Range("P4").Value = '11.4.1_S2'
$ws5.Range("Q4").Value = '1.2.1_S2'
$ws5.Range("R4").Value = '5.10.2_S2'
$ws5.Range("S4").Value = '1.6.1_S1'
$ws5.Range("T4").Value = '7.1.1_S2'
$ws5.Range("U4").Value = '7.1.1_S1'
$ws5.Range("V4").Value = '9.9.2_S1'
$ws5.Range("W4").Value = '4.5.1_S1'
$ws5.Range("X4").Value = 'rs'
$ws5.Range("Y4").Value = 'st'
$ws5.Range("Z4").Value = 'bcm'

# Row 5: counts per sample
$ws5.Range("C5").Value = 2
$ws5.Range("D5").Value = 2
$ws5.Range("E5").Value = 2
$ws5.Range("F5").Value = 2
$ws5.Range("G5").Value = 2
$ws5.Range("H5").Value = 2
$ws5.Range("I5").Value = 2
$ws5.Range("J5").Value = 2
$ws5.Range("K5").Value = 2
$ws5.Range("L5").Value = 2
$ws5.Range("M5").Value = 2
$ws5.Range("N5").Value = 4
$ws5.Range("O5").Value = 6
$ws5.Range("P5").Value = 7
$ws5.Range("Q5").Value = 9
$ws5.Range("R5").Value = 10
$ws5.Range("S5").Value = 11
$ws5.Range("T5").Value = 11
$ws5.Range("U5").Value = 12
$ws5.Range("V5").Value = 13
$ws5.Range("W5").Value = 14
$ws5.Range("X5").Value = 18
$ws5.Range("Y5").Value = 41
$ws5.Range("Z5").Value = 42

# Row 10: headers for samples with count > 10 (replicated across multiple IR files)
$ws5.Range("F10").Value = '1.6.1_S1'
$ws5.Range("G10").Value = '7.1.1_S2'
$ws5.Range("H10").Value = '7.1.1_S1'
$ws5.Range("I10").Value = '9.9.2_S1'
$ws5.Range("J10").Value = '4.5.1_S1'
$ws5.Range("K10").Value = 'rs'
$ws5.Range("L10").Value = 'st'
$ws5.Range("M10").Value = 'bcm'

# Rows 11-12: combined summary strings
$ws5.Range("D11").Value = '3.6.1_S1   6.1.1_S1   8.1.2_S1   8.3.1_S1   8.3.2_S1   8.5.2_S1   8.6.1_S1   8.6.2_S1   8.7.2_S1   8.8.1_S1   9.4.2_S1 11.10.2_S1 11.10.1_S1  11.4.1_S2   1.2.1_S2  5.10.2_S2 '
$ws5.Range("D12").Value = '         2          2          2          2          2          2          2          2          2          2          2          4          6          7          9         10'

# Row 18: table header for soil-sample replicate counts
$ws5.Range("B18").Value = 'Soil Samle'
$ws5.Range("C18").Value = 'n IR Files '

# Rows 19-42: soil sample name + replicate (IR file) count
$ws5.Range("B19").Value = '3.6.1_S1'
$ws5.Range("C19").Value = 2
$ws5.Range("B20").Value = '6.1.1_S1'
$ws5.Range("C20").Value = 2
$ws5.Range("B21").Value = '8.1.2_S1'
$ws5.Range("C21").Value = 2
$ws5.Range("B22").Value = '8.3.1_S1'
$ws5.Range("C22").Value = 2
$ws5.Range("B23").Value = '8.3.2_S1'
$ws5.Range("C23").Value = 2
$ws5.Range("B24").Value = '8.5.2_S1'
$ws5.Range("C24").Value = 2
$ws5.Range("B25").Value = '8.6.1_S1'
$ws5.Range("C25").Value = 2
$ws5.Range("B26").Value = '8.6.2_S1'
$ws5.Range("C26").Value = 2
$ws5.Range("B27").Value = '8.7.2_S1'
$ws5.Range("C27").Value = 2
$ws5.Range("B28").Value = '8.8.1_S1'
$ws5.Range("C28").Value = 2
$ws5.Range("B29").Value = '9.4.2_S1'
$ws5.Range("C29").Value = 2
$ws5.Range("B30").Value = '11.10.2_S1'
$ws5.Range("C30").Value = 4
$ws5.Range("B31").Value = '11.10.1_S1'
$ws5.Range("C31").Value = 6
$ws5.Range("B32").Value = '11.4.1_S2'
$ws5.Range("C32").Value = 7
$ws5.Range("B33").Value = '1.2.1_S2'
$ws5.Range("C33").Value = 9
$ws5.Range("B34").Value = '5.10.2_S2'
$ws5.Range("C34").Value = 10
$ws5.Range("B35").Value = '1.6.1_S1'
$ws5.Range("C35").Value = 11
$ws5.Range("B36").Value = '7.1.1_S2'
$ws5.Range("C36").Value = 11
$ws5.Range("B37").Value = '7.1.1_S1'
$ws5.Range("C37").Value = 12
$ws5.Range("B38").Value = '9.9.2_S1'
$ws5.Range("C38").Value = 13
$ws5.Range("B39").Value = '4.5.1_S1'
$ws5.Range("C39").Value = 14
$ws5.Range("B40").Value = 'rs'
$ws5.Range("C40").Value = 18
$ws5.Range("B41").Value = 'st'
$ws5.Range("C41").Value = 41
$ws5.Range("B42").Value = 'bcm'
$ws5.Range("C42").Value = 42

# Highlight the two outlier replicate counts (rows 30-31) in red
$ws5.Range("B30:C31").Font.Color = 255

# Column B width (auto-fit to sample-name text)
$ws5.Columns.Item(2).ColumnWidth = 9.15

# View: scroll / selection
$ws5.Activate()
$win5 = $excel.ActiveWindow
$win5.ScrollRow = 16
$win5.ScrollColumn = 1
$ws5.Range("B18:C42").Select()

# --- Make 'Summaries' the active tab and set its final selection ---
$ws4.Activate()
$ws4.Range("C19").Select()
